# Slovenia Prva Liga base update (10-06-2024 21:53)
# The source data feed re-fetched fixtures and, due to upstream ordering
# changes, rows 2/3 and rows 177/178 ended up swapped in the dataset.
# Apply the same swap here: exchange the full row content (columns B:AD)
# between row 2 <-> row 3 and between row 177 <-> row 178, leaving the
# running index in column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")

    $valsA = $rangeA.Value2()
    $valsB = $rangeB.Value2()

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-Rows $ws 2 3
Swap-Rows $ws 177 178
